$wb = $excel.ActiveWorkbook

# Move the "总计" worksheet to be before the "2021-Q3" worksheet,
# so the tab order becomes: 总计, 2021-Q3
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2021-Q3")

$totalSheet.Move($q3Sheet)
